# Added two new fields to NSDE: `inactivation_date` and `reactivation_date`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook / sheet is being repurposed for the NSDE field reference.
$ws.Name = "nsde_fields"

# --- Row 12: inactivation_date -----------------------------------------
$ws.Range("B12").Value = "inactivation_date"
$ws.Range("C12").Value = "string"
$ws.Range("D12").Value = "The date on which registration or listing data was inactivated by FDA due to inaccuracies, incompleteness or incompliance."

# B12 picks up a plain (non-wrapping, non-indented) style rather than the
# "left/top" style used by the other Field Name cells.
$ws.Range("B12").HorizontalAlignment = 1
$ws.Range("B12").VerticalAlignment = -4107

$ws.Rows.Item(12).RowHeight = 34

# --- Row 13: reactivation_date ------------------------------------------
$ws.Range("B13").Value = "reactivation_date"
$ws.Range("C13").Value = "string"
$ws.Range("D13").Value = "The date on which a previously FDA inactivated registration or listing data is reactivated."

$ws.Rows.Item(13).RowHeight = 17

# Leave the view parked on the newly added field, zoomed in a bit, matching
# the author's final selection/zoom when they saved.
[void]$ws.Range("B12").Select()
$excel.ActiveWindow.Zoom = 140
